# Adds the 13 missing daily COVID report rows (04/12/2020 .. 16/12/2020)
# to the bottom of the data table, extending it from row 250 to row 263.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each inner array: rowNumber, data, dia, casos, mortes, Ativos,
#                   taxa morte contaminados, Curados, Casos negativos,
#                   Testes realizados, novosCasos, suspeitos,
#                   mortesSuspeitas, suspeitosAtivos, novosTestes,
#                   leitos_clinicos_ocupados, leitos_uti_ocupados, semana
$data = @(
    ,@(251, "04/12/2020", 250, 3968, 100, 81, 0.0252016129032258, 3770, 12548, 16516, 44, 313, 1, 312, 164, 25, 5, 36)
    ,@(252, "05/12/2020", 251, 4007, 100, 58, 0.02495632642874969, 3832, 12608, 16615, 39, 248, 2, 246, 99, 19, 5, 36)
    ,@(253, "06/12/2020", 252, 4029, 100, 48, 0.02482005460412013, 3864, 12680, 16709, 22, 226, 4, 222, 94, 15, 4, 37)
    ,@(254, "07/12/2020", 253, 4049, 100, 35, 0.02469745616201531, 3897, 12680, 16729, 20, 234, 5, 229, 20, 15, 6, 37)
    ,@(255, "08/12/2020", 254, 4076, 100, 41, 0.02453385672227674, 3918, 12704, 16780, 27, 183, 6, 177, 51, 14, 6, 37)
    ,@(256, "09/12/2020", 255, 4097, 100, 46, 0.0244081034903588, 3934, 12801, 16898, 21, 276, 6, 270, 118, 18, 7, 37)
    ,@(257, "10/12/2020", 256, 4141, 100, 59, 0.02414875633904854, 3965, 12925, 17066, 44, 221, 7, 214, 168, 20, 10, 37)
    ,@(258, "11/12/2020", 257, 4180, 101, 66, 0.02416267942583732, 3996, 12963, 17143, 39, 314, 2, 312, 77, 18, 10, 37)
    ,@(259, "12/12/2020", 258, 4180, 101, 66, 0.02416267942583732, 3996, 12963, 17143, 0, 314, 2, 312, 0, 18, 10, 37)
    ,@(260, "13/12/2020", 259, 4185, 101, 33, 0.02413381123058542, 4034, 12974, 17159, 5, 316, 3, 313, 16, 19, 10, 38)
    ,@(261, "14/12/2020", 260, 4185, 101, 33, 0.02413381123058542, 4034, 12974, 17159, 0, 316, 3, 313, 0, 19, 10, 38)
    ,@(262, "15/12/2020", 261, 4228, 102, 36, 0.02412488174077578, 4073, 13076, 17304, 43, 446, 1, 445, 145, 21, 7, 38)
    ,@(263, "16/12/2020", 262, 4262, 102, 38, 0.02393242609103707, 4105, 13173, 17435, 34, 435, 2, 434, 131, 16, 7, 38)
)

foreach ($row in $data) {
    $r = $row[0]

    # Column A holds a date-shaped string (e.g. "04/12/2020") that must stay
    # literal text, not get auto-converted into a date serial number.
    $cellA = $ws.Cells.Item($r, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value = $row[1]
    $cellA.Style = "Normal"

    $ws.Cells.Item($r, 2).Value  = $row[2]
    $ws.Cells.Item($r, 3).Value  = $row[3]
    $ws.Cells.Item($r, 4).Value  = $row[4]
    $ws.Cells.Item($r, 5).Value  = $row[5]
    $ws.Cells.Item($r, 6).Value  = $row[6]
    $ws.Cells.Item($r, 7).Value  = $row[7]
    $ws.Cells.Item($r, 8).Value  = $row[8]
    $ws.Cells.Item($r, 9).Value  = $row[9]
    $ws.Cells.Item($r, 10).Value = $row[10]
    $ws.Cells.Item($r, 11).Value = $row[11]
    $ws.Cells.Item($r, 12).Value = $row[12]
    $ws.Cells.Item($r, 13).Value = $row[13]
    $ws.Cells.Item($r, 14).Value = $row[14]
    $ws.Cells.Item($r, 15).Value = $row[15]
    $ws.Cells.Item($r, 16).Value = $row[16]
    $ws.Cells.Item($r, 17).Value = $row[17]
}
